# Update column F (dSF) values for the mikolas_miles sheet based on repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 3
    3  = 1
    4  = -1
    5  = -7
    6  = 3
    7  = 1
    9  = 1
    10 = -2
    11 = -1
    12 = -1
    13 = 2
    14 = -2
    15 = 2
    16 = -1
    17 = -6
    18 = 1
    19 = 2
    21 = 2
    22 = -2
    23 = -1
    27 = -2
    28 = 3
    29 = -1
    30 = -2
    31 = -2
    32 = 1
    33 = -4
    34 = -1
    35 = 6
    36 = 2
    37 = 4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
